$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the single AnimCube import line with three new component
#    imports (Exhibit, YouTube, ImageCollage).
# ---------------------------------------------------------------------------
$importRange = $d.Content
$importRange.Find.ClearFormatting()
$foundImport = $importRange.Find.Execute("import AnimCube from ""@site/src/components/AnimCube"";")

if ($foundImport) {
    $importXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t>import Exhibit from "@site/src/components/Exhibit";</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t>import YouTube from "@site/src/components/YouTube";</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t>import ImageCollage from '@site/src/components/ImageCollage';</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
    $null = $importRange.InsertXML($importXml)
}

# ---------------------------------------------------------------------------
# 2. Replace the single <AnimCube params="..." .../> tag with the new
#    multi-line <Exhibit ... /> component block (6 paragraphs).
# ---------------------------------------------------------------------------
$exhibitRange = $d.Content
$exhibitRange.Find.ClearFormatting()
$foundExhibit = $exhibitRange.Find.Execute('<AnimCube params="config=../../ExhibitConfig.txt&facelets=dldlyldldwwwwwwdlwdbbdbbdlddggdggdggdddooooooddddrddrr" width="400px" height="400px" />')

if ($foundExhibit) {
    $exhibitXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t>&lt;Exhibit</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t xml:space="preserve">  stickering={{</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t xml:space="preserve">    solved: "U D F B L R </w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>FL DFL DL DBL BL DF DB BR DBR</w:t></w:r>
            <w:r><w:t>",</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t xml:space="preserve">    orientedWithoutPermutation: "UL UF UR UB </w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>FR DR</w:t></w:r>
            <w:r><w:t>"</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t xml:space="preserve">  }}</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
            <w:r><w:t>/&gt;</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
    $null = $exhibitRange.InsertXML($exhibitXml)
}

Write-Host "Import found/replaced:" $foundImport
Write-Host "Exhibit found/replaced:" $foundExhibit
